# Apply cryptos list update (commit: "Updated cryptos list on Wed May 22 19:58:47 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.714.34'
$ws.Range("E2").Value = '  +0.48%  '

# Row 3
$ws.Range("D3").Value = '3.740.26'
$ws.Range("E3").Value = '  +0.69%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").Value = "'" + '612.85'
$ws.Range("E5").Value = '  +0.56%  '

# Row 6
$ws.Range("D6").Value = "'" + '178.16'
$ws.Range("E6").Value = '  +1.73%  '

# Row 7
$ws.Range("D7").Value = '3.738.49'
$ws.Range("E7").Value = '  +0.72%  '

# Row 9
$ws.Range("D9").Value = "'" + '0.529'
$ws.Range("E9").Value = '  -1.63%  '

# Row 10
$ws.Range("E10").Value = '  +0.53%  '

# Row 11
$ws.Range("D11").Value = "'" + '6.59'
$ws.Range("E11").Value = '  +3.88%  '

# Row 12
$ws.Range("E12").Value = '  -2.88%  '

# Row 13
$ws.Range("D13").Value = "'" + '39.91'
$ws.Range("E13").Value = '  -1.37%  '

# Row 14
$ws.Range("D14").Value = "'" + '0.0000254'
$ws.Range("E14").Value = '  +0.56%  '

# Row 15
$ws.Range("D15").Value = '4.361.39'
$ws.Range("E15").Value = '  +0.63%  '

# Row 16
$ws.Range("D16").Value = '3.739.61'
$ws.Range("E16").Value = '  +0.67%  '

# Row 17
$ws.Range("D17").Value = '69.732.10'
$ws.Range("E17").Value = '  +0.31%  '

# Row 18
$ws.Range("E18").Value = '  -2.38%  '

# Row 19
$ws.Range("D19").Value = "'" + '7.45'
$ws.Range("E19").Value = '  -1.34%  '

# Row 20
$ws.Range("D20").Value = "'" + '501.46'
$ws.Range("E20").Value = '  -2.08%  '

# Row 21
$ws.Range("D21").Value = "'" + '16.32'
$ws.Range("E21").Value = '  -2.14%  '

# Row 22
$ws.Range("E22").Value = '  -3.10%  '

# Row 23
$ws.Range("E23").Value = '  -0.40%  '

# Row 24
$ws.Range("E24").Value = '  +8.90%  '

# Row 25
$ws.Range("D25").Value = "'" + '85.99'
$ws.Range("E25").Value = '  -1.64%  '

# Row 26
$ws.Range("D26").Value = "'" + '11.79'
$ws.Range("E26").Value = '  +7.88%  '

# Row 27
$ws.Range("D27").Value = "'" + '12.90'
$ws.Range("E27").Value = '  -3.06%  '

# Row 28
$ws.Range("E28").Value = '  +8.68%  '

# Row 29
$ws.Range("E29").Value = '  +0.37%  '

# Row 30
$ws.Range("D30").Value = "'" + '2.47'
$ws.Range("E30").Value = '  -0.93%  '

# Row 31
$ws.Range("D31").Value = "'" + '8.15'
$ws.Range("E31").Value = '  +3.88%  '

# Row 32
$ws.Range("E32").Value = '  +3.17%  '

# Row 33
$ws.Range("D33").Value = "'" + '30.37'
$ws.Range("E33").Value = '  -2.02%  '

# Row 34
$ws.Range("E34").Value = '  -1.04%  '

# Row 35
$ws.Range("D35").Value = "'" + '0.999'
$ws.Range("E35").Value = '  -0.08%  '

# Row 36
$ws.Range("E36").Value = '  +1.89%  '

# Row 37
$ws.Range("E37").Value = '  -0.68%  '

# Row 38
$ws.Range("D38").Value = "'" + '0.356'
$ws.Range("E38").Value = '  +5.98%  '

# Row 39
$ws.Range("D39").Value = "'" + '0.138'
$ws.Range("E39").Value = '  +4.23%  '

# Row 40
$ws.Range("D40").Value = "'" + '453.38'
$ws.Range("E40").Value = '  +8.35%  '

# Row 41
$ws.Range("D41").Value = "'" + '3.09'
$ws.Range("E41").Value = '  +14.21%  '

# Row 42
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = "'" + '2.08'
$ws.Range("E42").Value = '  -3.94%  '

# Row 43
$ws.Range("B43").Value = 'Arweave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D43").Value = "'" + '46.09'
$ws.Range("E43").Value = '  +4.41%  '

# Row 44
$ws.Range("D44").Value = "'" + '49.72'

# Row 45
$ws.Range("E45").Value = '  -2.39%  '

# Row 46
$ws.Range("D46").Value = '2.948.93'
$ws.Range("E46").Value = '  -4.27%  '

# Row 47
$ws.Range("D47").Value = "'" + '0.0359'
$ws.Range("E47").Value = '  -0.72%  '

# Row 48
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = "'" + '1.00'
$ws.Range("E48").Value = '  -0.01%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = "'" + '27.12'
$ws.Range("E49").Value = '  -2.02%  '

# Row 50
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = "'" + '138.25'
$ws.Range("E50").Value = '  +2.76%  '

# Row 51
$ws.Range("D51").Value = "'" + '2.49'
$ws.Range("E51").Value = '  -0.13%  '
